$wb = $excel.ActiveWorkbook

# Sheet 1 updates
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 42
$ws.Range("F5").Value = 172
$ws.Range("F6").Value = 1047
$ws.Range("F7").Value = 632
$ws.Range("F8").Value = 7976
$ws.Range("F10").Value = 192
$ws.Range("F11").Value = 6834
$ws.Range("F12").Value = 160
$ws.Range("F13").Value = 298
$ws.Range("F14").Value = 4877
$ws.Range("F17").Value = 5291
$ws.Range("F19").Value = 315
$ws.Range("F20").Value = 316
$ws.Range("F21").Value = 435
$ws.Range("F26").Value = 94
$ws.Range("F27").Value = 8987
$ws.Range("F29").Value = 1604
$ws.Range("F31").Value = 39
$ws.Range("F33").Value = 825
$ws.Range("F35").Value = 70
$ws.Range("F36").Value = 1004
$ws.Range("F37").Value = 1146
$ws.Range("F38").Value = 50
$ws.Range("F39").Value = 4685
$ws.Range("F42").Value = 1153
$ws.Range("F44").Value = 141
$ws.Range("F45").Value = 70
$ws.Range("F47").Value = 1232
$ws.Range("F48").Value = 28
$ws.Range("F49").Value = 57

# Sheet 2 updates
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 37
$ws.Range("F9").Value = 177
$ws.Range("F17").Value = 885

# Sheet 4 updates
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 42
$ws.Range("F6").Value = 172
$ws.Range("F7").Value = 37
$ws.Range("F8").Value = 1047
$ws.Range("F9").Value = 632
$ws.Range("F10").Value = 7976
$ws.Range("F12").Value = 192
$ws.Range("F13").Value = 6834
$ws.Range("F14").Value = 160
$ws.Range("F15").Value = 298
$ws.Range("F17").Value = 4878
$ws.Range("F19").Value = 5292
$ws.Range("F21").Value = 315
$ws.Range("F22").Value = 316
$ws.Range("F23").Value = 435
$ws.Range("F28").Value = 94
$ws.Range("F29").Value = 177
$ws.Range("F30").Value = 8987
$ws.Range("F32").Value = 1604
$ws.Range("F33").Value = 39
$ws.Range("F35").Value = 825
$ws.Range("F37").Value = 70
$ws.Range("F38").Value = 1004
$ws.Range("F39").Value = 1146
$ws.Range("F40").Value = 50
$ws.Range("F41").Value = 4685
$ws.Range("F43").Value = 1153
$ws.Range("F44").Value = 141
$ws.Range("F45").Value = 70
$ws.Range("F47").Value = 1232
$ws.Range("F48").Value = 28
$ws.Range("F49").Value = 57
